$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 19:20:29"
$ws1.Range("A3").Value = "Total filas: 564"

$s1data = @(
    @("19:20:18", "19:22", "16_SANTA ANA", 2, "LP1912", "30/12/2025"),
    @("19:20:18", "19:22", "26_HERNANDEZ", 2, "LP1912", "30/12/2025"),
    @("19:20:18", "19:29", "15_ABASTO", 9, "LP1912", "30/12/2025"),
    @("19:20:18", "19:31", "16_SANTA ANA", 11, "LP1912", "30/12/2025"),
    @("19:20:18", "19:40", "215C_EL PATO", 20, "LP1912", "30/12/2025"),
    @("19:20:18", "19:41", "14_ABASTO", 21, "LP1912", "30/12/2025"),
    @("19:20:18", "19:41", "16_SANTA ANA", 21, "LP1912", "30/12/2025"),
    @("19:20:18", "19:51", "11X44_ETCHEVERRY", 31, "LP1912", "30/12/2025"),
    @("19:20:18", "19:51", "16_P MOR-SANTA ANA", 31, "LP1912", "30/12/2025"),
    @("19:20:18", "19:52", "81_EL PELIGRO", 32, "LP1912", "30/12/2025"),
    @("19:20:18", "20:00", "17_ROMERO", 40, "LP1912", "30/12/2025"),
    @("19:20:18", "20:01", "14_ABASTO", 41, "LP1912", "30/12/2025"),
    @("19:20:18", "20:08", "10_OLMOS", 48, "LP1912", "30/12/2025"),
    @("19:20:18", "20:10", "15_ABASTO", 50, "LP1912", "30/12/2025"),
    @("19:20:18", "20:11", "16_P MOR-167 Y 521", 51, "LP1912", "30/12/2025"),
    @("19:20:18", "20:13", "23_HERNANDEZ", 53, "LP1912", "30/12/2025"),
    @("19:20:18", "20:21", "26_HERNANDEZ", 61, "LP1912", "30/12/2025"),
    @("19:20:18", "20:23", "11_ETCHEVERRY", 63, "LP1912", "30/12/2025"),
    @("19:20:18", "20:24", "215A_EL PATO", 64, "LP1912", "30/12/2025"),
    @("19:20:18", "20:53", "15_ABASTO", 93, "LP1912", "30/12/2025"),
    @("19:20:18", "20:56", "10_OLMOS", 96, "LP1912", "30/12/2025"),
    @("19:20:18", "20:57", "23_HERNANDEZ", 97, "LP1912", "30/12/2025"),
)
$r = 544
foreach ($row in $s1data) {
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Value = $row[1]
    $ws1.Cells.Item($r, 4).Value = $row[2]
    $ws1.Cells.Item($r, 5).Value = $row[3]
    $ws1.Cells.Item($r, 6).Value = $row[4]
    $ws1.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 19:20:29"
$ws2.Range("A3").Value = "Total filas: 37"

$s2data = @(
    @("30/12/2025", "19:20:18", "19:40", "215C_EL PATO", 20, "LP1912"),
    @("30/12/2025", "19:20:18", "20:24", "215A_EL PATO", 64, "LP1912"),
)
$r = 37
foreach ($row in $s2data) {
    $ws2.Cells.Item($r, 2).Value = $row[0]
    $ws2.Cells.Item($r, 3).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
    $ws2.Cells.Item($r, 6).Value = $row[4]
    $ws2.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 30/12/2025 19:20:29"
$ws3.Range("A3").Value = "Total filas: 71"

$s3data = @(
    @("30/12/2025", "19:20:23", "19:54", "215C_LA PLATA", 34, "L6203"),
)
$r = 72
foreach ($row in $s3data) {
    $ws3.Cells.Item($r, 2).Value = $row[0]
    $ws3.Cells.Item($r, 3).Value = $row[1]
    $ws3.Cells.Item($r, 4).Value = $row[2]
    $ws3.Cells.Item($r, 5).Value = $row[3]
    $ws3.Cells.Item($r, 6).Value = $row[4]
    $ws3.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
